$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Absent" column (H) to be the complement of the "Real" column (E):
# when a student was really present (E=1) Absent=0, otherwise Absent=1.
# This forms the consolidated report by filling in/correcting the Absent values.
for ($row = 3; $row -le 21; $row++) {
    $real = $ws.Cells.Item($row, 5).Value2
    $ws.Cells.Item($row, 8).Value2 = 1 - $real
}
